$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2,1).Value = "ECs"
$ws.Cells.Item(2,2).Value = "Edn1"
$ws.Cells.Item(2,3).Value = "Ednra"
$ws.Cells.Item(2,4).Value = "ECs"
$ws.Cells.Item(2,5).Value = 3
$ws.Cells.Item(2,6).Value = 1
$ws.Cells.Item(2,7).Value = 9.957023333333334
$ws.Cells.Item(2,8).Value = 29.87107
$ws.Cells.Item(2,9).Value = 0.7839926662698464
$ws.Cells.Item(2,10).Value = 0.7839926662698464
$ws.Cells.Item(2,11).Value = 3
$ws.Cells.Item(2,12).Value = 1
$ws.Cells.Item(2,13).Value = 1.508394
$ws.Cells.Item(2,14).Value = 4.525182
$ws.Cells.Item(2,15).Value = 0.0276475339394655
$ws.Cells.Item(2,16).Value = 0.0276475339394655
$ws.Cells.Item(2,17).Value = 15.01911425386
$ws.Cells.Item(2,18).Value = 135.17202828474
$ws.Cells.Item(2,19).Value = 0.02167546384898763
$ws.Cells.Item(2,20).Value = 0.02167546384898763

# Row 3
$ws.Cells.Item(3,1).Value = "ECs"
$ws.Cells.Item(3,2).Value = "Edn1"
$ws.Cells.Item(3,3).Value = "Ednra"
$ws.Cells.Item(3,4).Value = "FAPs"
$ws.Cells.Item(3,5).Value = 3
$ws.Cells.Item(3,6).Value = 1
$ws.Cells.Item(3,7).Value = 9.957023333333334
$ws.Cells.Item(3,8).Value = 29.87107
$ws.Cells.Item(3,9).Value = 0.7839926662698464
$ws.Cells.Item(3,10).Value = 0.7839926662698464
$ws.Cells.Item(3,11).Value = 3
$ws.Cells.Item(3,12).Value = 1
$ws.Cells.Item(3,13).Value = 8.961352
$ws.Cells.Item(3,14).Value = 26.884056
$ws.Cells.Item(3,15).Value = 0.1642536920482958
$ws.Cells.Item(3,16).Value = 0.1642536920482958
$ws.Cells.Item(3,17).Value = 89.22839096221334
$ws.Cells.Item(3,18).Value = 803.05551865992
$ws.Cells.Item(3,19).Value = 0.1287736899736097
$ws.Cells.Item(3,20).Value = 0.1287736899736097

# Row 4
$ws.Cells.Item(4,1).Value = "ECs"
$ws.Cells.Item(4,2).Value = "Edn1"
$ws.Cells.Item(4,3).Value = "Ednra"
$ws.Cells.Item(4,4).Value = "sCs"
$ws.Cells.Item(4,5).Value = 3
$ws.Cells.Item(4,6).Value = 1
$ws.Cells.Item(4,7).Value = 9.957023333333334
$ws.Cells.Item(4,8).Value = 29.87107
$ws.Cells.Item(4,9).Value = 0.7839926662698464
$ws.Cells.Item(4,10).Value = 0.7839926662698464
$ws.Cells.Item(4,11).Value = 3
$ws.Cells.Item(4,12).Value = 1
$ws.Cells.Item(4,13).Value = 44.08824833333333
$ws.Cells.Item(4,14).Value = 132.264745
$ws.Cells.Item(4,15).Value = 0.8080987740122387
$ws.Cells.Item(4,16).Value = 0.8080987740122386
$ws.Cells.Item(4,17).Value = 438.9877173807944
$ws.Cells.Item(4,18).Value = 3950.88945642715
$ws.Cells.Item(4,19).Value = 0.6335435124472492
$ws.Cells.Item(4,20).Value = 0.633543512447249

# Row 5
$ws.Cells.Item(5,1).Value = "FAPs"
$ws.Cells.Item(5,2).Value = "Edn1"
$ws.Cells.Item(5,3).Value = "Ednra"
$ws.Cells.Item(5,4).Value = "ECs"
$ws.Cells.Item(5,5).Value = 3
$ws.Cells.Item(5,6).Value = 1
$ws.Cells.Item(5,7).Value = 2.539481333333333
$ws.Cells.Item(5,8).Value = 7.618444
$ws.Cells.Item(5,9).Value = 0.1999528046497
$ws.Cells.Item(5,10).Value = 0.1999528046497
$ws.Cells.Item(5,11).Value = 3
$ws.Cells.Item(5,12).Value = 1
$ws.Cells.Item(5,13).Value = 1.508394
$ws.Cells.Item(5,14).Value = 4.525182
$ws.Cells.Item(5,15).Value = 0.0276475339394655
$ws.Cells.Item(5,16).Value = 0.0276475339394655
$ws.Cells.Item(5,17).Value = 3.830538406312
$ws.Cells.Item(5,18).Value = 34.474845656808
$ws.Cells.Item(5,19).Value = 0.005528201952843896
$ws.Cells.Item(5,20).Value = 0.005528201952843896

# Row 6
$ws.Cells.Item(6,1).Value = "FAPs"
$ws.Cells.Item(6,2).Value = "Edn1"
$ws.Cells.Item(6,3).Value = "Ednra"
$ws.Cells.Item(6,4).Value = "FAPs"
$ws.Cells.Item(6,5).Value = 3
$ws.Cells.Item(6,6).Value = 1
$ws.Cells.Item(6,7).Value = 2.539481333333333
$ws.Cells.Item(6,8).Value = 7.618444
$ws.Cells.Item(6,9).Value = 0.1999528046497
$ws.Cells.Item(6,10).Value = 0.1999528046497
$ws.Cells.Item(6,11).Value = 3
$ws.Cells.Item(6,12).Value = 1
$ws.Cells.Item(6,13).Value = 8.961352
$ws.Cells.Item(6,14).Value = 26.884056
$ws.Cells.Item(6,15).Value = 0.1642536920482958
$ws.Cells.Item(6,16).Value = 0.1642536920482958
$ws.Cells.Item(6,17).Value = 22.75718612542933
$ws.Cells.Item(6,18).Value = 204.814675128864
$ws.Cells.Item(6,19).Value = 0.03284298639912486
$ws.Cells.Item(6,20).Value = 0.03284298639912486

# Row 7
$ws.Cells.Item(7,1).Value = "FAPs"
$ws.Cells.Item(7,2).Value = "Edn1"
$ws.Cells.Item(7,3).Value = "Ednra"
$ws.Cells.Item(7,4).Value = "sCs"
$ws.Cells.Item(7,5).Value = 3
$ws.Cells.Item(7,6).Value = 1
$ws.Cells.Item(7,7).Value = 2.539481333333333
$ws.Cells.Item(7,8).Value = 7.618444
$ws.Cells.Item(7,9).Value = 0.1999528046497
$ws.Cells.Item(7,10).Value = 0.1999528046497
$ws.Cells.Item(7,11).Value = 3
$ws.Cells.Item(7,12).Value = 1
$ws.Cells.Item(7,13).Value = 44.08824833333333
$ws.Cells.Item(7,14).Value = 132.264745
$ws.Cells.Item(7,15).Value = 0.8080987740122387
$ws.Cells.Item(7,16).Value = 0.8080987740122386
$ws.Cells.Item(7,17).Value = 111.9612836618644
$ws.Cells.Item(7,18).Value = 1007.65155295678
$ws.Cells.Item(7,19).Value = 0.1615816162977312
$ws.Cells.Item(7,20).Value = 0.1615816162977312

# Row 8
$ws.Cells.Item(8,1).Value = "sCs"
$ws.Cells.Item(8,2).Value = "Edn1"
$ws.Cells.Item(8,3).Value = "Ednra"
$ws.Cells.Item(8,4).Value = "ECs"
$ws.Cells.Item(8,5).Value = 2
$ws.Cells.Item(8,6).Value = 0.6666666666666666
$ws.Cells.Item(8,7).Value = 0.203899
$ws.Cells.Item(8,8).Value = 0.611697
$ws.Cells.Item(8,9).Value = 0.01605452908045364
$ws.Cells.Item(8,10).Value = 0.01605452908045364
$ws.Cells.Item(8,11).Value = 3
$ws.Cells.Item(8,12).Value = 1
$ws.Cells.Item(8,13).Value = 1.508394
$ws.Cells.Item(8,14).Value = 4.525182
$ws.Cells.Item(8,15).Value = 0.0276475339394655
$ws.Cells.Item(8,16).Value = 0.0276475339394655
$ws.Cells.Item(8,17).Value = 0.307560028206
$ws.Cells.Item(8,18).Value = 2.768040253854
$ws.Cells.Item(8,19).Value = 0.0004438681376339779
$ws.Cells.Item(8,20).Value = 0.0004438681376339778

# Row 9
$ws.Cells.Item(9,1).Value = "sCs"
$ws.Cells.Item(9,2).Value = "Edn1"
$ws.Cells.Item(9,3).Value = "Ednra"
$ws.Cells.Item(9,4).Value = "FAPs"
$ws.Cells.Item(9,5).Value = 2
$ws.Cells.Item(9,6).Value = 0.6666666666666666
$ws.Cells.Item(9,7).Value = 0.203899
$ws.Cells.Item(9,8).Value = 0.611697
$ws.Cells.Item(9,9).Value = 0.01605452908045364
$ws.Cells.Item(9,10).Value = 0.01605452908045364
$ws.Cells.Item(9,11).Value = 3
$ws.Cells.Item(9,12).Value = 1
$ws.Cells.Item(9,13).Value = 8.961352
$ws.Cells.Item(9,14).Value = 26.884056
$ws.Cells.Item(9,15).Value = 0.1642536920482958
$ws.Cells.Item(9,16).Value = 0.1642536920482958
$ws.Cells.Item(9,17).Value = 1.827210711448
$ws.Cells.Item(9,18).Value = 16.444896403032
$ws.Cells.Item(9,19).Value = 0.002637015675561241
$ws.Cells.Item(9,20).Value = 0.002637015675561241

# Row 10
$ws.Cells.Item(10,1).Value = "sCs"
$ws.Cells.Item(10,2).Value = "Edn1"
$ws.Cells.Item(10,3).Value = "Ednra"
$ws.Cells.Item(10,4).Value = "sCs"
$ws.Cells.Item(10,5).Value = 2
$ws.Cells.Item(10,6).Value = 0.6666666666666666
$ws.Cells.Item(10,7).Value = 0.203899
$ws.Cells.Item(10,8).Value = 0.611697
$ws.Cells.Item(10,9).Value = 0.01605452908045364
$ws.Cells.Item(10,10).Value = 0.01605452908045364
$ws.Cells.Item(10,11).Value = 3
$ws.Cells.Item(10,12).Value = 1
$ws.Cells.Item(10,13).Value = 44.08824833333333
$ws.Cells.Item(10,14).Value = 132.264745
$ws.Cells.Item(10,15).Value = 0.8080987740122387
$ws.Cells.Item(10,16).Value = 0.8080987740122386
$ws.Cells.Item(10,17).Value = 8.989549746918334
$ws.Cells.Item(10,18).Value = 80.90594772226501
$ws.Cells.Item(10,19).Value = 0.01297364526725842
$ws.Cells.Item(10,20).Value = 0.01297364526725842

